## Cargo.xlsx maintenance edit:
##  - rename the "order" column header to "Order"
##  - convert the repeated =A2..=A7 formulas in column F into one shared formula
##  - move the active selection to B2 (saved view state)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header B1 used to read "order" (lower-case); the data now groups it with
# the other proper-case headers (Name, CanRaid, IsArmy, Pic) as "Order".
$ws.Range("B1").Value = "Order"

# F2:F7 all just mirror column A ("=A2", "=A3", ... "=A7"); re-enter them as
# a single range formula so Excel stores them as one shared formula group.
$ws.Range("F2:F7").Formula = "=A2"

# Restore the saved cursor position to B2 for this sheet view.
$ws.Range("B2").Select() | Out-Null
